# Regenerate orders with updated distance/size codes.
# Distance codes: D64 -> D69, D80 -> D86, D51 -> D55
# Size code:      S30 -> S31
# These substrings appear inside Condition, Filename_Left, Filename_Right,
# Distance and Size columns (and nowhere else), so a straightforward
# find/replace over the sheet's used range reproduces the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange

# Order matters only in that each pair of old/new tokens must not collide
# with another rule's replacement text; these four are mutually disjoint.
$used.Replace("D64", "D69")
$used.Replace("D80", "D86")
$used.Replace("D51", "D55")
$used.Replace("S30", "S31")
